$d = $word.ActiveDocument

# 1. Update "Total de citas programadas: 4" -> "...: 5"
#    Scope the Find to just the paragraph that holds this label so the
#    bold label run and the plain count run stay separate (matches the
#    original two-run layout).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Total de citas programadas:*") {
        $p.Range.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "5", 2)
    }
}

# 2. Rework the appointments table.
$t = $d.Tables.Item(1)

# Row 2 (08:30 - 08:45): comprador COLFRESH COFFEE -> REGIONAL S.A.S
$t.Rows.Item(2).Cells.Item(3).Range.Text = "REGIONAL S.A.S"

# Row 3: time 08:45 - 09:00 -> 09:00 - 09:15 ; comprador PROCOLOMBIA -> COLFRESH COFFEE
$t.Rows.Item(3).Cells.Item(1).Range.Text = "09:00 - 09:15"
$t.Rows.Item(3).Cells.Item(3).Range.Text = "COLFRESH COFFEE"

# Row 4: time 09:00 - 09:15 -> 09:30 - 09:45 ; comprador CAFÉ MOLINA stays the same
$t.Rows.Item(4).Cells.Item(1).Range.Text = "09:30 - 09:45"

# Row 5: time 09:30 - 09:45 -> 10:00 - 10:15 ; comprador REGIONAL S.A.S -> PROCOLOMBIA
$t.Rows.Item(5).Cells.Item(1).Range.Text = "10:00 - 10:15"
$t.Rows.Item(5).Cells.Item(3).Range.Text = "PROCOLOMBIA"

# New row 6: 10:15 - 10:30 | (blank mesa) | BOX BRAND
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "10:15 - 10:30"
$newRow.Cells.Item(3).Range.Text = "BOX BRAND"
